$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the crypto price/volume snapshot (and a handful of rank swaps)
# to match the latest scrape. A leading apostrophe is used for the few
# numeric-looking price strings that Excel would otherwise normalize
# (stripping trailing zeros / switching to scientific notation), so the
# literal text is preserved exactly like the source cells.

$ws.Range("D2").Value = "26.168.73"
$ws.Range("E2").Value = "  -1.18%  "

$ws.Range("D3").Value = "1.810.60"
$ws.Range("E3").Value = "  -1.53%  "

$ws.Range("D4").Value = "1.007"
$ws.Range("E4").Value = "  +0.64%  "

$ws.Range("D5").Value = "240.96"
$ws.Range("E5").Value = "  -7.41%  "

$ws.Range("D6").Value = "1.011"
$ws.Range("E6").Value = "  +0.97%  "

$ws.Range("D7").Value = "0.5102"
$ws.Range("E7").Value = "  -2.75%  "

$ws.Range("D8").Value = "0.2463"
$ws.Range("E8").Value = "  -22.87%  "

$ws.Range("D9").Value = "0.06147"
$ws.Range("E9").Value = "  -9.39%  "

$ws.Range("D10").Value = "1.846.13"
$ws.Range("E10").Value = "  +0.68%  "

$ws.Range("D11").Value = "0.06862"
$ws.Range("E11").Value = "  -11.45%  "

$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").Value = "14.81"
$ws.Range("E12").Value = "  -20.99%  "

$ws.Range("B13").Value = "Litecoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D13").Value = "79.23"
$ws.Range("E13").Value = "  -9.71%  "

$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "0.5984"
$ws.Range("E14").Value = "  -23.77%  "

$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "4.404"
$ws.Range("E15").Value = "  -12.14%  "

$ws.Range("D16").Value = "'1.010"
$ws.Range("E16").Value = "  +0.91%  "

$ws.Range("D17").Value = "1.011"
$ws.Range("E17").Value = "  +0.98%  "

$ws.Range("D18").Value = "26.190.45"
$ws.Range("E18").Value = "  -1.19%  "

$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").Value = "11.37"
$ws.Range("E19").Value = "  -17.83%  "

$ws.Range("B20").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C20").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D20").Value = "2.052.96"
$ws.Range("E20").Value = "  -0.91%  "

$ws.Range("D21").Value = "'0.000005901"
$ws.Range("E21").Value = "  -25.73%  "

$ws.Range("D22").Value = "3.952"
$ws.Range("E22").Value = "  -14.55%  "

$ws.Range("D23").Value = "5.208"
$ws.Range("E23").Value = "  -12.68%  "

$ws.Range("D24").Value = "7.948"
$ws.Range("E24").Value = "  -15.10%  "

$ws.Range("D25").Value = "131.74"
$ws.Range("E25").Value = "  -7.01%  "

$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").Value = "1.876"
$ws.Range("E26").Value = "  -13.76%  "

$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "14.53"
$ws.Range("E27").Value = "  -14.13%  "

$ws.Range("B28").Value = "BitcoinCash"
$ws.Range("C28").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D28").Value = "99.17"
$ws.Range("E28").Value = "  -11.03%  "

$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "1.245"
$ws.Range("E29").Value = "  -26.21%  "

$ws.Range("D30").Value = "0.08346"
$ws.Range("E30").Value = "  -3.96%  "

$ws.Range("D31").Value = "3.634"
$ws.Range("E31").Value = "  -12.64%  "

$ws.Range("E32").Value = "  -3.36%  "

$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "3.169"
$ws.Range("E33").Value = "  -22.17%  "

$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "0.04273"
$ws.Range("E34").Value = "  -12.45%  "

$ws.Range("D35").Value = "1.058"
$ws.Range("E35").Value = "  -6.70%  "

$ws.Range("D36").Value = "2.958"
$ws.Range("E36").Value = "  -4.39%  "

$ws.Range("D37").Value = "0.6251"
$ws.Range("E37").Value = "  -14.20%  "

$ws.Range("D38").Value = "2.087"
$ws.Range("E38").Value = "  -6.85%  "

$ws.Range("D39").Value = "1.011"
$ws.Range("E39").Value = "  +1.00%  "

$ws.Range("D40").Value = "0.8379"
$ws.Range("E40").Value = "  -6.15%  "

$ws.Range("D41").Value = "99.89"
$ws.Range("E41").Value = "  -8.85%  "

$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "0.01459"
$ws.Range("E42").Value = "  -16.75%  "

$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "5.191"
$ws.Range("E43").Value = "  -12.44%  "

$ws.Range("D44").Value = "'0.3830"
$ws.Range("E44").Value = "  -19.76%  "

$ws.Range("D45").Value = "0.05271"
$ws.Range("E45").Value = "  -9.87%  "

$ws.Range("D46").Value = "'6.180"
$ws.Range("E46").Value = "  -19.35%  "

$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "53.42"
$ws.Range("E47").Value = "  -10.23%  "

$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").Value = "29.86"
$ws.Range("E48").Value = "  -14.32%  "

$ws.Range("B49").Value = "USDD"
$ws.Range("C49").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D49").Value = "1.011"
$ws.Range("E49").Value = "  +0.74%  "

$ws.Range("D50").Value = "1.008"
$ws.Range("E50").Value = "  +0.77%  "

$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "0.1032"
$ws.Range("E51").Value = "  -16.13%  "
